# Append two new high-value dataset rows to Sheet1 and refresh the view.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 90 - MGNREGA (Gram Panchayat / Yearly)
$ws.Range("A90").Value = "mord-mnrega-gp-yr-abc"
$ws.Range("B90").Value = "Mahatma Gandhi National Rural Employment Guarantee Act (MGNREGA) "

# New row 91 - PMAY-G (Gram Panchayat / Other)
$ws.Range("B91").Value = "Pradhan Mantri Awaas Yojana - Gramin"
$ws.Range("A91").Value = "mohua-pmayg-gp-ol-abc"

# Remaining columns (granularity_level, frequency) for both new rows
$ws.Range("C90").Value = "Gram Panchayat"
$ws.Range("D90").Value = "Yearly"
$ws.Range("C91").Value = "Gram Panchayat"
$ws.Range("D91").Value = "Other"

# Match the yellow highlight styling used for other manually-added rows
$ws.Range("A90:D91").Interior.Color = 65535

# Remove the stale autofilter now that new unfiltered rows were appended
$ws.AutoFilterMode = $false | Out-Null

# Leave selection on the last entered cell, as in the source edit
$ws.Range("E91").Select() | Out-Null
